# Applies the diff: update row 73 (E73/F73), append rows 74-82 to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing timestamp cell style (custom date format) from F73
$tsFormat = $ws.Range("F73").NumberFormat

# --- Row 73: E73 becomes a plain phone number; F73 timestamp refined ---
$ws.Range("E73").Value = 919510038048
$ws.Range("F73").Value = 45987.63049465277

# --- Row 74 ---
$ws.Range("A74").Value = "[Call Started]"
$ws.Range("D74").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("E74").Value = 919510038048
$ws.Range("F74").Value = 45988.66858461806
$ws.Range("F74").NumberFormat = $tsFormat

# --- Row 75 ---
$ws.Range("A75").Value = "[Call Started]"
$ws.Range("D75").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("E75").Value = 919510038048
$ws.Range("F75").Value = 45988.6697115625
$ws.Range("F75").NumberFormat = $tsFormat

# --- Row 76 ---
$ws.Range("A76").Value = "[Intro response]"
$ws.Range("B76").Value = "Yash"
$ws.Range("C76").Value = "neutral"
$ws.Range("D76").Value = "Here are our latest offers:
- Laptop Pro 
- Smart watch 
- Bluetooth Earbuds
Which product would you like to purchase?"
$ws.Range("E76").Value = 919510038048
$ws.Range("F76").Value = 45988.67003817129
$ws.Range("F76").NumberFormat = $tsFormat

# --- Row 77 ---
$ws.Range("A77").Value = "[Fallback]"
$ws.Range("B77").Value = "Bluetooth earphones"
$ws.Range("C77").Value = "neutral"
$ws.Range("D77").Value = "Sorry, we don’t have that product right now.
Here are our latest offers:
- Laptop Pro : one of the best laptop you can get right now with high end specs at ₹75000
- Smart watch : Best watch in market with all your daily tracking at ₹12000
- Bluetooth Earbuds: best anc earbuds with this price point  at ₹4000
Which product would you like to purchase?"
$ws.Range("E77").Value = 919510038048
$ws.Range("F77").Value = 45988.67015928241
$ws.Range("F77").NumberFormat = $tsFormat

# --- Row 78 ---
$ws.Range("A78").Value = "[Fallback]"
$ws.Range("B78").Value = "smart watch"
$ws.Range("C78").Value = "neutral"
$ws.Range("D78").Value = "Sorry, we don’t have that product right now.
Here are our latest offers:
- Laptop Pro : one of the best laptop you can get right now with high end specs at ₹75000
- Smart watch : Best watch in market with all your daily tracking at ₹12000
- Bluetooth Earbuds: best anc earbuds with this price point  at ₹4000
Which product would you like to purchase?"
$ws.Range("E78").Value = 919510038048
$ws.Range("F78").Value = 45988.67049449074
$ws.Range("F78").NumberFormat = $tsFormat

# --- Row 79 ---
$ws.Range("A79").Value = "[Fallback]"
$ws.Range("B79").Value = "smart watch"
$ws.Range("C79").Value = "neutral"
$ws.Range("D79").Value = "Sorry, we don’t have that product right now.
Here are our latest offers:
- Laptop Pro : one of the best laptop you can get right now with high end specs at ₹75000
- Smart watch : Best watch in market with all your daily tracking at ₹12000
- Bluetooth Earbuds: best anc earbuds with this price point  at ₹4000
Which product would you like to purchase?"
$ws.Range("E79").Value = 919510038048
$ws.Range("F79").Value = 45988.67072487268
$ws.Range("F79").NumberFormat = $tsFormat

# --- Row 80 ---
$ws.Range("A80").Value = "[Call Started]"
$ws.Range("D80").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("E80").Value = 919510038048
$ws.Range("F80").Value = 45988.67656203704
$ws.Range("F80").NumberFormat = $tsFormat

# --- Row 81 ---
$ws.Range("A81").Value = "[Intro response]"
$ws.Range("B81").Value = "Yash"
$ws.Range("C81").Value = "neutral"
$ws.Range("D81").Value = "Here are our latest offers:
- Laptop Pro 
- Smart watch 
- Bluetooth Earbuds
Which product would you like to purchase?"
$ws.Range("E81").Value = 919510038048
$ws.Range("F81").Value = 45988.67678699074
$ws.Range("F81").NumberFormat = $tsFormat

# --- Row 82 ---
$ws.Range("A82").Value = "[Product match]"
$ws.Range("B82").Value = "laptop"
$ws.Range("C82").Value = "neutral"
$ws.Range("D82").Value = "Great choice! I’ve sent the link of Laptop Pro  to your phone number ending with 8048. Thank you for your time! I really appreciate it."
$ws.Range("E82").Value = "'+919510038048"
$ws.Range("E82").Style = "Normal"
$ws.Range("F82").Value = 45988.67693401646
$ws.Range("F82").NumberFormat = $tsFormat
